$p = $ppt.ActivePresentation

# --- Slide 35 ("How to write computer algorithms") ---
$s = $p.Slides.Item(35)

# 1) Touch the speaker notes for this slide (mirrors the author opening the
#    Notes pane for this slide - PowerPoint mints a new, still-empty
#    notesSlide part + relationship as a side effect).
$notesBody = $s.NotesPage.Shapes.AddPlaceholder(2)

# 2) Edit the body textbox ("TextBox 6"): split out "syntax" as italic and
#    change "...understand the semantics Python..." to
#    "...understand the semantics of Python...".
$sh = $s.Shapes.Item(6)
$tr = $sh.TextFrame.TextRange

# " what is and isn't valid Python syntax. A second main job..." ->
# italicise the standalone word "syntax" (2nd occurrence in the text).
$tr.Characters(201, 6).Font.Italic = $true

# "...understand the semantics Python i.e..." -> insert "of " right after
# "semantics " (still italic, as it lands inside that italic run).
[void]$tr.Characters(248, 10).InsertAfter("of ")
